$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 46048.01041666666
$ws.Cells.Item(2, 2).Value = 0
$ws.Cells.Item(3, 1).Value = 46048.02083333334
$ws.Cells.Item(3, 2).Value = 0
$ws.Cells.Item(4, 1).Value = 46048.03125
$ws.Cells.Item(4, 2).Value = 0
$ws.Cells.Item(5, 1).Value = 46048.04166666666
$ws.Cells.Item(5, 2).Value = 0
$ws.Cells.Item(6, 1).Value = 46048.05208333334
$ws.Cells.Item(6, 2).Value = 0
$ws.Cells.Item(7, 1).Value = 46048.0625
$ws.Cells.Item(7, 2).Value = 0
$ws.Cells.Item(8, 1).Value = 46048.07291666666
$ws.Cells.Item(8, 2).Value = 0
$ws.Cells.Item(9, 1).Value = 46048.08333333334
$ws.Cells.Item(9, 2).Value = 0
$ws.Cells.Item(10, 1).Value = 46048.09375
$ws.Cells.Item(10, 2).Value = 0
$ws.Cells.Item(11, 1).Value = 46048.10416666666
$ws.Cells.Item(11, 2).Value = 0
$ws.Cells.Item(12, 1).Value = 46048.11458333334
$ws.Cells.Item(12, 2).Value = 0
$ws.Cells.Item(13, 1).Value = 46048.125
$ws.Cells.Item(13, 2).Value = 0
$ws.Cells.Item(14, 1).Value = 46048.13541666666
$ws.Cells.Item(14, 2).Value = 0
$ws.Cells.Item(15, 1).Value = 46048.14583333334
$ws.Cells.Item(15, 2).Value = 0
$ws.Cells.Item(16, 1).Value = 46048.15625
$ws.Cells.Item(16, 2).Value = 0
$ws.Cells.Item(17, 1).Value = 46048.16666666666
$ws.Cells.Item(17, 2).Value = 0
$ws.Cells.Item(18, 1).Value = 46048.17708333334
$ws.Cells.Item(18, 2).Value = 0
$ws.Cells.Item(19, 1).Value = 46048.1875
$ws.Cells.Item(19, 2).Value = 0
$ws.Cells.Item(20, 1).Value = 46048.19791666666
$ws.Cells.Item(20, 2).Value = 0
$ws.Cells.Item(21, 1).Value = 46048.20833333334
$ws.Cells.Item(21, 2).Value = 0
$ws.Cells.Item(22, 1).Value = 46048.21875
$ws.Cells.Item(22, 2).Value = 0.536
$ws.Cells.Item(23, 1).Value = 46048.22916666666
$ws.Cells.Item(23, 2).Value = 0.5570000000000001
$ws.Cells.Item(24, 1).Value = 46048.23958333334
$ws.Cells.Item(24, 2).Value = 0.582
$ws.Cells.Item(25, 1).Value = 46048.25
$ws.Cells.Item(25, 2).Value = 0.639
$ws.Cells.Item(26, 1).Value = 46048.26041666666
$ws.Cells.Item(26, 2).Value = 1.442
$ws.Cells.Item(27, 1).Value = 46048.27083333334
$ws.Cells.Item(27, 2).Value = 1.867
$ws.Cells.Item(28, 1).Value = 46048.28125
$ws.Cells.Item(28, 2).Value = 3.464
$ws.Cells.Item(29, 1).Value = 46048.29166666666
$ws.Cells.Item(29, 2).Value = 4.644
$ws.Cells.Item(30, 1).Value = 46048.30208333334
$ws.Cells.Item(30, 2).Value = 11.934
$ws.Cells.Item(31, 1).Value = 46048.3125
$ws.Cells.Item(31, 2).Value = 21.887
$ws.Cells.Item(32, 1).Value = 46048.32291666666
$ws.Cells.Item(32, 2).Value = 32.596
$ws.Cells.Item(33, 1).Value = 46048.33333333334
$ws.Cells.Item(33, 2).Value = 47.499
$ws.Cells.Item(34, 1).Value = 46048.34375
$ws.Cells.Item(34, 2).Value = 79.43300000000001
$ws.Cells.Item(35, 1).Value = 46048.35416666666
$ws.Cells.Item(35, 2).Value = 99.657
$ws.Cells.Item(36, 1).Value = 46048.36458333334
$ws.Cells.Item(36, 2).Value = 118.805
$ws.Cells.Item(37, 1).Value = 46048.375
$ws.Cells.Item(37, 2).Value = 135.737
$ws.Cells.Item(38, 1).Value = 46048.38541666666
$ws.Cells.Item(38, 2).Value = 171.616
$ws.Cells.Item(39, 1).Value = 46048.39583333334
$ws.Cells.Item(39, 2).Value = 189.065
$ws.Cells.Item(40, 1).Value = 46048.40625
$ws.Cells.Item(40, 2).Value = 207.018
$ws.Cells.Item(41, 1).Value = 46048.41666666666
$ws.Cells.Item(41, 2).Value = 222.369
$ws.Cells.Item(42, 1).Value = 46048.42708333334
$ws.Cells.Item(42, 2).Value = 246.333
$ws.Cells.Item(43, 1).Value = 46048.4375
$ws.Cells.Item(43, 2).Value = 259.098
$ws.Cells.Item(44, 1).Value = 46048.44791666666
$ws.Cells.Item(44, 2).Value = 271.294
$ws.Cells.Item(45, 1).Value = 46048.45833333334
$ws.Cells.Item(45, 2).Value = 281.039
$ws.Cells.Item(46, 1).Value = 46048.46875
$ws.Cells.Item(46, 2).Value = 282.114
$ws.Cells.Item(47, 1).Value = 46048.47916666666
$ws.Cells.Item(47, 2).Value = 287.303
$ws.Cells.Item(48, 1).Value = 46048.48958333334
$ws.Cells.Item(48, 2).Value = 289.749
$ws.Cells.Item(49, 1).Value = 46048.5
$ws.Cells.Item(49, 2).Value = 287.808
$ws.Cells.Item(50, 1).Value = 46048.51041666666
$ws.Cells.Item(50, 2).Value = 284.657
$ws.Cells.Item(51, 1).Value = 46048.52083333334
$ws.Cells.Item(51, 2).Value = 275.379
$ws.Cells.Item(52, 1).Value = 46048.53125
$ws.Cells.Item(52, 2).Value = 265.319
$ws.Cells.Item(53, 1).Value = 46048.54166666666
$ws.Cells.Item(53, 2).Value = 253.471
$ws.Cells.Item(54, 1).Value = 46048.55208333334
$ws.Cells.Item(54, 2).Value = 236.085
$ws.Cells.Item(55, 1).Value = 46048.5625
$ws.Cells.Item(55, 2).Value = 219.937
$ws.Cells.Item(56, 1).Value = 46048.57291666666
$ws.Cells.Item(56, 2).Value = 204.218
$ws.Cells.Item(57, 1).Value = 46048.58333333334
$ws.Cells.Item(57, 2).Value = 182.366
$ws.Cells.Item(58, 1).Value = 46048.59375
$ws.Cells.Item(58, 2).Value = 153.488
$ws.Cells.Item(59, 1).Value = 46048.60416666666
$ws.Cells.Item(59, 2).Value = 133.769
$ws.Cells.Item(60, 1).Value = 46048.61458333334
$ws.Cells.Item(60, 2).Value = 112.921
$ws.Cells.Item(61, 1).Value = 46048.625
$ws.Cells.Item(61, 2).Value = 94.44799999999999
$ws.Cells.Item(62, 1).Value = 46048.63541666666
$ws.Cells.Item(62, 2).Value = 60.614
$ws.Cells.Item(63, 1).Value = 46048.64583333334
$ws.Cells.Item(63, 2).Value = 47.74
$ws.Cells.Item(64, 1).Value = 46048.65625
$ws.Cells.Item(64, 2).Value = 37.29
$ws.Cells.Item(65, 1).Value = 46048.66666666666
$ws.Cells.Item(65, 2).Value = 27.387
$ws.Cells.Item(66, 1).Value = 46048.67708333334
$ws.Cells.Item(66, 2).Value = 19.445
$ws.Cells.Item(67, 1).Value = 46048.6875
$ws.Cells.Item(67, 2).Value = 8.84
$ws.Cells.Item(68, 1).Value = 46048.69791666666
$ws.Cells.Item(68, 2).Value = 7.532
$ws.Cells.Item(69, 1).Value = 46048.70833333334
$ws.Cells.Item(69, 2).Value = 6.648
$ws.Cells.Item(70, 1).Value = 46048.71875
$ws.Cells.Item(70, 2).Value = 0.68
$ws.Cells.Item(71, 1).Value = 46048.72916666666
$ws.Cells.Item(71, 2).Value = 0.669
$ws.Cells.Item(72, 1).Value = 46048.73958333334
$ws.Cells.Item(72, 2).Value = 0.662
$ws.Cells.Item(73, 1).Value = 46048.75
$ws.Cells.Item(73, 2).Value = 0.671
$ws.Cells.Item(74, 1).Value = 46048.76041666666
$ws.Cells.Item(74, 2).Value = 0.65
$ws.Cells.Item(75, 1).Value = 46048.77083333334
$ws.Cells.Item(75, 2).Value = 0
$ws.Cells.Item(76, 1).Value = 46048.78125
$ws.Cells.Item(76, 2).Value = 0
$ws.Cells.Item(77, 1).Value = 46048.79166666666
$ws.Cells.Item(77, 2).Value = 0
$ws.Cells.Item(78, 1).Value = 46048.80208333334
$ws.Cells.Item(78, 2).Value = 0
$ws.Cells.Item(79, 1).Value = 46048.8125
$ws.Cells.Item(79, 2).Value = 0
$ws.Cells.Item(80, 1).Value = 46048.82291666666
$ws.Cells.Item(80, 2).Value = 0
$ws.Cells.Item(81, 1).Value = 46048.83333333334
$ws.Cells.Item(81, 2).Value = 0
$ws.Cells.Item(82, 1).Value = 46048.84375
$ws.Cells.Item(82, 2).Value = 0
$ws.Cells.Item(83, 1).Value = 46048.85416666666
$ws.Cells.Item(83, 2).Value = 0
$ws.Cells.Item(84, 1).Value = 46048.86458333334
$ws.Cells.Item(84, 2).Value = 0
$ws.Cells.Item(85, 1).Value = 46048.875
$ws.Cells.Item(85, 2).Value = 0
$ws.Cells.Item(86, 1).Value = 46048.88541666666
$ws.Cells.Item(86, 2).Value = 0.45
$ws.Cells.Item(87, 1).Value = 46048.89583333334
$ws.Cells.Item(87, 2).Value = 0
$ws.Cells.Item(88, 1).Value = 46048.90625
$ws.Cells.Item(88, 2).Value = 0
$ws.Cells.Item(89, 1).Value = 46048.91666666666
$ws.Cells.Item(89, 2).Value = 0
$ws.Cells.Item(90, 1).Value = 46048.92708333334
$ws.Cells.Item(90, 2).Value = 0
$ws.Cells.Item(91, 1).Value = 46048.9375
$ws.Cells.Item(91, 2).Value = 0
$ws.Cells.Item(92, 1).Value = 46048.94791666666
$ws.Cells.Item(92, 2).Value = 0
$ws.Cells.Item(93, 1).Value = 46048.95833333334
$ws.Cells.Item(93, 2).Value = 0
$ws.Cells.Item(94, 1).Value = 46048.96875
$ws.Cells.Item(94, 2).Value = 0
$ws.Cells.Item(95, 1).Value = 46048.97916666666
$ws.Cells.Item(95, 2).Value = 0
$ws.Cells.Item(96, 1).Value = 46048.98958333334
$ws.Cells.Item(96, 2).Value = 0
$ws.Cells.Item(97, 1).Value = 46049
$ws.Cells.Item(97, 2).Value = 0
